# =========================================================================
# Applies the "add ground_warfae files and data" edit to the workbook.
#  - ship sheet: add column F (values mirroring column E pattern), tweak a
#    number of existing values, and remove the green highlight fill that
#    used to cover C:E (now also not applied to F).
#  - SSM sheet: update a few numeric values.
#  - SAM sheet: update a few numeric values.
#  - inception sheet: update a value and widen column A.
#  - selections (active cell) updated on several sheets.
# =========================================================================

$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------
# Sheet "ship"
# -------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("ship")
$ws1.Activate()

# --- value corrections on existing columns (B:E) ---
$ws1.Range("B7").Value2  = 36
$ws1.Range("C7").Value2  = 36
$ws1.Range("D7").Value2  = 36
$ws1.Range("E7").Value2  = 36

$ws1.Range("B10").Value2 = 15
$ws1.Range("E10").Value2 = 20

$ws1.Range("B12").Value2 = 15

$ws1.Range("C13").Value2 = 1
$ws1.Range("D13").Value2 = 1
$ws1.Range("E13").Value2 = 2

$ws1.Range("B14").Value2 = 1
$ws1.Range("C14").Value2 = 5
$ws1.Range("D14").Value2 = 5
$ws1.Range("E14").Value2 = 5

$ws1.Range("B24").Value2 = 150
$ws1.Range("C24").Value2 = 150
$ws1.Range("D24").Value2 = 150
$ws1.Range("E24").Value2 = 150

$ws1.Range("B34").Value2 = 2
$ws1.Range("C34").Value2 = 2
$ws1.Range("D34").Value2 = 2
$ws1.Range("E34").Value2 = 2

$ws1.Range("B35").Value2 = 3
$ws1.Range("C35").Value2 = 3
$ws1.Range("D35").Value2 = 3
$ws1.Range("E35").Value2 = 3

# --- new column F ---
$ws1.Range("F1").Value2  = 5
$ws1.Range("F2").Value2  = 25
$ws1.Range("F3").Value2  = 90
$ws1.Range("F4").Value2  = 3
$ws1.Range("F5").Value2  = 10
$ws1.Range("F6").Value2  = 24
$ws1.Range("F7").Value2  = 36
$ws1.Range("F8").Value2  = 20
$ws1.Range("F9").Value2  = 3
$ws1.Range("F10").Value2 = 20
$ws1.Range("F11").Value2 = 1
$ws1.Range("F12").Value2 = 20
$ws1.Range("F13").Value2 = 2
$ws1.Range("F14").Value2 = 5
$ws1.Range("F15").Value2 = 8
$ws1.Range("F16").Value2 = 2
$ws1.Range("F17").Value2 = 20
$ws1.Range("F18").Value2 = -0.1
$ws1.Range("F19").Value2 = 10
$ws1.Range("F20").Value2 = 0.9
$ws1.Range("F21").Value2 = 100
$ws1.Range("F22").Value2 = 25
$ws1.Range("F23").Value2 = 20
$ws1.Range("F24").Value2 = 150
$ws1.Range("F25").Value2 = 6000000
$ws1.Range("F26").Value2 = 9300
$ws1.Range("F27").Value2 = 7.5
$ws1.Range("F28").Value2 = 4
$ws1.Range("F29").Value2 = 4
$ws1.Range("F30").Value2 = 4200
$ws1.Range("F31").Value2 = 1100
$ws1.Range("F32").Value2 = 10
$ws1.Range("F33").Value2 = 11
$ws1.Range("F34").Value2 = 2
$ws1.Range("F35").Value2 = 3
$ws1.Range("F36").Value2 = 2
$ws1.Range("F37").Value2 = 3
$ws1.Range("F38").Value2 = "yellow"

# --- remove the green fill / style (was style index 1, fill FF00B050) from
#     the whole C:F block, restoring default (unstyled) cells / columns ---
$ws1.Range("C1:F38").Style = "표준"
$ws1.Columns("C:F").ClearFormats()

# --- update the active cell / selection ---
$ws1.Range("F1").Select()

# -------------------------------------------------------------------
# Sheet "SSM"
# -------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("SSM")
$ws2.Activate()

$ws2.Range("B2").Value2 = 2.5
$ws2.Range("C2").Value2 = 3
$ws2.Range("D2").Value2 = 4.5

$ws2.Range("B8").Value2 = 150
$ws2.Range("C8").Value2 = 150
$ws2.Range("D8").Value2 = 150

$ws2.Range("C3").Select()

# -------------------------------------------------------------------
# Sheet "SAM"
# -------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("SAM")
$ws3.Activate()

$ws3.Range("B2").Value2 = 5.5
$ws3.Range("C2").Value2 = 5.5
$ws3.Range("D2").Value2 = 5.5
$ws3.Range("E2").Value2 = 5.5

$ws3.Range("B3").Value2 = 60
$ws3.Range("C3").Value2 = 60
$ws3.Range("D3").Value2 = 30
$ws3.Range("E3").Value2 = 30

$ws3.Range("N11").Select()

# -------------------------------------------------------------------
# Sheet "inception"
# -------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("inception")
$ws5.Activate()

$ws5.Range("B2").Value2 = 90

# widen column A (renders as width="16" in the saved XML)
$ws5.Columns("A").ColumnWidth = 15.29

$ws5.Range("B2").Select()

# -------------------------------------------------------------------
# Leave "ship" as the active/selected sheet, matching original tab state
# -------------------------------------------------------------------
$ws1.Activate()

Write-Output "edit complete"
